$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 36: new "Check BBR Values" field (FMR.BBR) -----------------------
# Shared-string insertion order (matches target: E=90 FMR.BBR, B=91 Check
# BBR Values, D=92 rich-text comment) is driven by write order below.
$ws.Range("E36").Value() = "FMR.BBR"
$ws.Range("B36").Value() = "Check BBR Values"

$commentIntro = 'Read text files in "_WorkflowReports_\texts" to evaluate the success of Boundary-Based Registration (FMR-VMR Coregistration) for each run. The results are displayed in the Command Window as well as written to the log file. BBR was the default option in BV20, but BV21 now defaults to use the NGF method instead so the new default is FALSE. '
$commentBold = 'Set TRUE if you used BBR for FMR-VMR Coregistration.'
$ws.Range("D36").Value() = ($commentIntro + $commentBold)
$boldChars = $ws.Range("D36").Characters($commentIntro.Length + 1, $commentBold.Length)
$boldChars.Font.Bold = $true

$ws.Range("C36").Value() = $false

# Row 36 grows to fit the long wrapped comment text.
$ws.Rows.Item(36).RowHeight = 90

# --- View state: scrolled down one field, new active cell -----------------
$ws.Activate()
$ws.Range("B36").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
